$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("6. évf")
$ws.Activate()

# New row 28: "do jigsaw a puzzle" / "kirakóssal játszik"
$ws.Cells.Item(28, 1).Value = "3a"
$ws.Cells.Item(28, 2).Value = "do jigsaw a puzzle"
$ws.Cells.Item(28, 3).Value = "kirakóssal játszik"

# Fix the mistranslated header in row 2 ("answer to phone" -> "answer the phone")
$ws.Range("B2").Value = "answer the phone"

# Append remaining new vocabulary rows (29-35)
$newRows = @(
    @("do karate", "karatézik"),
    @("go bowling", "bowlingozni megy"),
    @("go go-kart racing", "go-kart versenyre megy"),
    @("go ice skating", "korcsolyázik"),
    @("let's have fun", "Szórakozzunk/ érezzük jól magunkat"),
    @("indoor", "beltéri"),
    @("outdoor", "kültéri")
)

$r = 29
foreach ($pair in $newRows) {
    $ws.Cells.Item($r, 1).Value = "3a"
    $ws.Cells.Item($r, 2).Value = $pair[0]
    $ws.Cells.Item($r, 3).Value = $pair[1]
    $r = $r + 1
}

$ws.Range("A36").Select()
